$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 8847.5
$ws.Range("I74").Value = 8847.5
$ws.Range("K74").Value = 8847.5
$ws.Range("M74").Value = -7911.5

$ws.Range("H77").Value = 8847.5
$ws.Range("I77").Value = 8847.5
$ws.Range("K77").Value = 44237.5
$ws.Range("M77").Value = -39557.5

$ws.Range("H116").Value = 3785.5715
$ws.Range("I116").Value = 3833.3333
$ws.Range("K116").Value = 3833.3333
$ws.Range("M116").Value = -391.3332999999998

$ws.Range("H132").Value = 1959.0588
$ws.Range("I132").Value = 1593.8572
$ws.Range("K132").Value = 4781.571599999999
$ws.Range("M132").Value = -2251.571599999999

$ws.Range("H137").Value = 5657.231
$ws.Range("I137").Value = 5994.6665
$ws.Range("K137").Value = 17983.9995
$ws.Range("M137").Value = -15433.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3382.077
$ws.Range("I61").Value = 2997.25
$ws.Range("K61").Value = 2997.25
$ws.Range("M61").Value = -2785.25

$ws.Range("H102").Value = 464.57144
$ws.Range("I102").Value = 445.05264
$ws.Range("K102").Value = 445.05264
$ws.Range("M102").Value = 1176.94736

$ws.Range("H132").Value = 3059.0293
$ws.Range("I132").Value = 3103.3333
$ws.Range("J132").Value = 2726.75
$ws.Range("K132").Value = 9309.999899999999
$ws.Range("L132").Value = 8180.25
$ws.Range("M132").Value = -6779.999899999999
$ws.Range("N132").Value = -13240.25

$ws.Range("H136").Value = 3382.077
$ws.Range("I136").Value = 2997.25
$ws.Range("K136").Value = 8991.75
$ws.Range("M136").Value = -6441.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1567.3334
$ws.Range("I99").Value = 1182.5
$ws.Range("K99").Value = 1182.5
$ws.Range("M99").Value = 315.5

$ws.Range("H107").Value = 4700.75
$ws.Range("I107").Value = 1895
$ws.Range("J107").Value = 7506.5
$ws.Range("K107").Value = 1895
$ws.Range("L107").Value = 7506.5
$ws.Range("M107").Value = 25
$ws.Range("N107").Value = -11346.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2248.3
$ws.Range("I31").Value = 2081.3333
$ws.Range("J31").Value = 2498.75
$ws.Range("K31").Value = 2081.3333
$ws.Range("L31").Value = 2498.75
$ws.Range("M31").Value = -1786.3333
$ws.Range("N31").Value = -3088.75

$ws.Range("H34").Value = 2248.3
$ws.Range("I34").Value = 2081.3333
$ws.Range("J34").Value = 2498.75
$ws.Range("K34").Value = 2081.3333
$ws.Range("L34").Value = 2498.75
$ws.Range("M34").Value = -1879.3333
$ws.Range("N34").Value = -2902.75

$ws.Range("H58").Value = 2493.8
$ws.Range("I58").Value = 1957.6
$ws.Range("J58").Value = 4102.4
$ws.Range("K58").Value = 1957.6
$ws.Range("L58").Value = 4102.4
$ws.Range("M58").Value = -1754.6
$ws.Range("N58").Value = -4508.4

$ws.Range("H86").Value = 7497.625
$ws.Range("I86").Value = 7996.4
$ws.Range("K86").Value = 7996.4
$ws.Range("M86").Value = -6873.4

$ws.Range("H89").Value = 7497.625
$ws.Range("I89").Value = 7996.4
$ws.Range("K89").Value = 39982
$ws.Range("M89").Value = -34366

$ws.Range("H94").Value = 1493.3334
$ws.Range("I94").Value = 1590.25
$ws.Range("J94").Value = 1299.5
$ws.Range("K94").Value = 1590.25
$ws.Range("L94").Value = 1299.5
$ws.Range("M94").Value = -1139.25
$ws.Range("N94").Value = -2201.5

$ws.Range("H122").Value = 5033
$ws.Range("I122").Value = 5646.8
$ws.Range("K122").Value = 16940.4
$ws.Range("M122").Value = -14490.4

$ws.Range("H132").Value = 1816.3334
$ws.Range("I132").Value = 1816.3334
$ws.Range("K132").Value = 5449.0002
$ws.Range("M132").Value = -2919.0002

$ws.Range("H134").Value = 2540.2104
$ws.Range("I134").Value = 2575.2222
$ws.Range("K134").Value = 7725.6666
$ws.Range("M134").Value = -5190.6666

$ws.Range("H136").Value = 2493.8
$ws.Range("I136").Value = 1957.6
$ws.Range("J136").Value = 4102.4
$ws.Range("K136").Value = 5872.799999999999
$ws.Range("L136").Value = 12307.2
$ws.Range("M136").Value = -3322.799999999999
$ws.Range("N136").Value = -17407.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 105.25
$ws.Range("I10").Value = 105.25
$ws.Range("K10").Value = 315.75
$ws.Range("M10").Value = -176.75

$ws.Range("H11").Value = 248.66667
$ws.Range("I11").Value = 248.66667
$ws.Range("K11").Value = 746.00001
$ws.Range("M11").Value = -606.00001

$ws.Range("H38").Value = 263.8
$ws.Range("I38").Value = 263.8
$ws.Range("K38").Value = 791.4000000000001
$ws.Range("M38").Value = -444.4000000000001

$ws.Range("H39").Value = 8833.333000000001
$ws.Range("J39").Value = 8833.333000000001
$ws.Range("L39").Value = 26499.999
$ws.Range("N39").Value = -27087.999

$ws.Range("H40").Value = 23.142857
$ws.Range("I40").Value = 26.666666
$ws.Range("J40").Value = 2
$ws.Range("K40").Value = 106.666664
$ws.Range("L40").Value = 8
$ws.Range("M40").Value = -37.666664
$ws.Range("N40").Value = -146

$ws.Range("H47").Value = 252.2
$ws.Range("I47").Value = 252.2
$ws.Range("K47").Value = 756.5999999999999
$ws.Range("M47").Value = -325.5999999999999

$ws.Range("H50").Value = 319.2857
$ws.Range("I50").Value = 280.83334
$ws.Range("K50").Value = 842.5000200000001
$ws.Range("M50").Value = -361.5000200000001

$ws.Range("H53").Value = 319.2857
$ws.Range("I53").Value = 280.83334
$ws.Range("K53").Value = 842.5000200000001
$ws.Range("M53").Value = -361.5000200000001

$ws.Range("H116").Value = 68119.8
$ws.Range("I116").Value = 68119.8
$ws.Range("K116").Value = 204359.4
$ws.Range("M116").Value = -200917.4

$ws.Range("H117").Value = 422.83334
$ws.Range("J117").Value = 422
$ws.Range("L117").Value = 1266
$ws.Range("N117").Value = -8150

$ws.Range("H122").Value = 750.6
$ws.Range("I122").Value = 686.6
$ws.Range("J122").Value = 814.6
$ws.Range("K122").Value = 6179.400000000001
$ws.Range("L122").Value = 7331.400000000001
$ws.Range("M122").Value = -3729.400000000001
$ws.Range("N122").Value = -12231.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null

$ws.Range("H132").Value = 1748.1333
$ws.Range("I132").Value = 1906.25
$ws.Range("J132").Value = 1115.6666
$ws.Range("K132").Value = 5718.75
$ws.Range("L132").Value = 3346.9998
$ws.Range("M132").Value = -3188.75
$ws.Range("N132").Value = -8406.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = $null

$ws.Range("H46").Value = 858.1818
$ws.Range("I46").Value = 860.2222
$ws.Range("K46").Value = 860.2222
$ws.Range("M46").Value = -672.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 41353.5
$ws.Range("J45").Value = 56370.75
$ws.Range("L45").Value = 56370.75
$ws.Range("N45").Value = -57352.75

$ws.Range("H81").Value = 3171.2083
$ws.Range("J81").Value = 2399
$ws.Range("L81").Value = 4798
$ws.Range("N81").Value = -6920

$ws.Range("H84").Value = 3171.2083
$ws.Range("J84").Value = 2399
$ws.Range("L84").Value = 23990
$ws.Range("N84").Value = -34598

$ws.Range("H102").Value = 33333
$ws.Range("J102").Value = 33333
$ws.Range("L102").Value = 33333
$ws.Range("N102").Value = -39823

$ws.Range("H107").Value = 613.7143
$ws.Range("J107").Value = 533.6
$ws.Range("L107").Value = 1600.8
$ws.Range("N107").Value = -5440.8

$ws.Range("H132").Value = 9005.177
$ws.Range("J132").Value = 19477
$ws.Range("L132").Value = 58431
$ws.Range("N132").Value = -63491
